$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Header 2 (header1.xml) - BTec logo: docPr/cNvPr id="1" - rename image1.jpg -> image2.jpg
$hdr2 = $sec.Headers.Item(2)
$hShape = $hdr2.Range.InlineShapes.Item(1)
$hShape.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.jpg"

# Footer 2 (footer1.xml) - Pearson logo: docPr/cNvPr id="3" - rename image2.png -> image1.png
$ftr2 = $sec.Footers.Item(2)
$fShape2 = $ftr2.Range.InlineShapes.Item(1)
$fShape2.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.png"

# Footer 1 (footer2.xml) - Pearson logo: docPr/cNvPr id="2" - rename image2.png -> image1.png
$ftr1 = $sec.Footers.Item(1)
$fShape1 = $ftr1.Range.InlineShapes.Item(1)
$fShape1.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.png"
